# Populate Sheet1 with the flight-search test-data table (header row + one
# data row) that the commit added, reproducing xl/sharedStrings.xml's
# first-seen string order along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (left to right; "Mobile_no" and "from" each repeat an
# already-used string so no new shared-string entries are created for them).
$headers = @("from","where","Mobile_no","Mobile_no","error_message","from","to","adults","children","infants","travel_class")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - sample data. Columns are written in A,B,C,D,F,G,H,I,J,K,E order
# (E2 last) so the shared-string table ends up in the exact sequence the
# workbook ships with: "Please enter a valid number" is appended after
# Chennai/Mumbai/Economy rather than right after Manali.
$ws.Cells.Item(2, 1).Value = "Bengaluru"
$ws.Cells.Item(2, 2).Value = "Manali"
$ws.Cells.Item(2, 3).Value = 8015993932
$ws.Cells.Item(2, 4).Value = 12345
$ws.Cells.Item(2, 6).Value = "Chennai"
$ws.Cells.Item(2, 7).Value = "Mumbai"
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = "Economy"
$ws.Cells.Item(2, 5).Value = "Please enter a valid number"

# Column C (Mobile_no) was manually narrowed to a width of 11 characters.
# ColumnWidth is stored with an ~0.8333 char padding offset added on save,
# so back it out here to land exactly on width="11" in the saved XML.
$ws.Columns.Item(3).ColumnWidth = 10.1666666666667

# Matches the saved <selection .../> range from the source workbook.
$ws.Range("L1:L4").Select()
